$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 281
$ws.Cells.Item(281, 6).Value = 80906
$ws.Cells.Item(281, 7).Value = 45988
$ws.Cells.Item(281, 8).Value = "E.03.0023"
$ws.Cells.Item(281, 9).Value = "JOELHEIRA ASSENTA,EMTO PISO"
$ws.Cells.Item(281, 10).Value = "UN"
$ws.Cells.Item(281, 11).Value = 10
$ws.Cells.Item(281, 12).Value = 79
$ws.Cells.Item(281, 13).Value = 790
$ws.Cells.Item(281, 14).Value = "00000000007459"
$ws.Cells.Item(281, 15).Value = "GADE FERRAMENTAS"

# Row 282
$ws.Cells.Item(282, 6).Value = 80905
$ws.Cells.Item(282, 7).Value = 45988
$ws.Cells.Item(282, 8).Value = "E.04.0073"
$ws.Cells.Item(282, 9).Value = "MARTELO DE UNHA PROFISSIONAL  COM CABO DE ALMA DE AÇO 25MM"
$ws.Cells.Item(282, 10).Value = "UN"
$ws.Cells.Item(282, 11).Value = 4
$ws.Cells.Item(282, 12).Value = 52.2
$ws.Cells.Item(282, 13).Value = 208.8
$ws.Cells.Item(282, 14).Value = "00000000000315"
$ws.Cells.Item(282, 15).Value = "CRUZADA"

# Row 283
$ws.Cells.Item(283, 6).Value = 80902
$ws.Cells.Item(283, 7).Value = 45988
$ws.Cells.Item(283, 8).Value = "E.04.0269"
$ws.Cells.Item(283, 9).Value = "BLOCO DE ESPUMA  OU CAMURÇA"
$ws.Cells.Item(283, 10).Value = "UN"
$ws.Cells.Item(283, 11).Value = 15
$ws.Cells.Item(283, 12).Value = 6
$ws.Cells.Item(283, 13).Value = 90
$ws.Cells.Item(283, 14).Value = "00000000008655"
$ws.Cells.Item(283, 15).Value = "DMC MATERIAIS"

# Row 284
$ws.Cells.Item(284, 6).Value = 80902
$ws.Cells.Item(284, 7).Value = 45988
$ws.Cells.Item(284, 8).Value = "E.04.1170"
$ws.Cells.Item(284, 9).Value = "LÁPIS PARA CARPINTEIRO IRWIN"
$ws.Cells.Item(284, 10).Value = "UN"
$ws.Cells.Item(284, 11).Value = 20
$ws.Cells.Item(284, 12).Value = 2.2
$ws.Cells.Item(284, 13).Value = 44
$ws.Cells.Item(284, 14).Value = "00000000008655"
$ws.Cells.Item(284, 15).Value = "DMC MATERIAIS"

# Row 285
$ws.Cells.Item(285, 6).Value = 80905
$ws.Cells.Item(285, 7).Value = 45988
$ws.Cells.Item(285, 8).Value = "J.06.0011"
$ws.Cells.Item(285, 9).Value = "REJUNTE ACRILICO"
$ws.Cells.Item(285, 10).Value = "KG"
$ws.Cells.Item(285, 11).Value = 10
$ws.Cells.Item(285, 12).Value = 26.6
$ws.Cells.Item(285, 13).Value = 266
$ws.Cells.Item(285, 14).Value = "00000000000315"
$ws.Cells.Item(285, 15).Value = "CRUZADA"

# Row 286
$ws.Cells.Item(286, 6).Value = 80905
$ws.Cells.Item(286, 7).Value = 45988
$ws.Cells.Item(286, 8).Value = "J.08.0002"
$ws.Cells.Item(286, 9).Value = "GESSO - 40 KG"
$ws.Cells.Item(286, 10).Value = "SC"
$ws.Cells.Item(286, 11).Value = 1
$ws.Cells.Item(286, 12).Value = 45
$ws.Cells.Item(286, 13).Value = 45
$ws.Cells.Item(286, 14).Value = "00000000000315"
$ws.Cells.Item(286, 15).Value = "CRUZADA"

# Row 287
$ws.Cells.Item(287, 6).Value = 80902
$ws.Cells.Item(287, 7).Value = 45988
$ws.Cells.Item(287, 8).Value = "S.10.0062"
$ws.Cells.Item(287, 9).Value = "PU 40  FLEX ADESIVO DE POLIURETANO DE CURA RÁPIDA COR BRANCO  TUBO 310ML"
$ws.Cells.Item(287, 10).Value = "UN"
$ws.Cells.Item(287, 11).Value = 50
$ws.Cells.Item(287, 12).Value = 17.62
$ws.Cells.Item(287, 13).Value = 881
$ws.Cells.Item(287, 14).Value = "00000000008655"
$ws.Cells.Item(287, 15).Value = "DMC MATERIAIS"

# Row 288
$ws.Cells.Item(288, 6).Value = 80905
$ws.Cells.Item(288, 7).Value = 45988
$ws.Cells.Item(288, 8).Value = "S.08.0200"
$ws.Cells.Item(288, 9).Value = "VIAPLUS 1000/TOP IMPER. BI-COMPONENTE(A+B)  - EMB. 18KG"
$ws.Cells.Item(288, 10).Value = "CX"
$ws.Cells.Item(288, 11).Value = 10
$ws.Cells.Item(288, 12).Value = 60.5
$ws.Cells.Item(288, 13).Value = 605
$ws.Cells.Item(288, 14).Value = "00000000000315"
$ws.Cells.Item(288, 15).Value = "CRUZADA"

# Row 289
$ws.Cells.Item(289, 6).Value = 80917
$ws.Cells.Item(289, 7).Value = 45988
$ws.Cells.Item(289, 8).Value = "S.08.0601"
$ws.Cells.Item(289, 9).Value = "BIANCO - ADESIVO DE ALTO DESEMPENHO PARA ARGAMASSAS E CHAPISCO - BALDE 18KG"
$ws.Cells.Item(289, 10).Value = "UN"
$ws.Cells.Item(289, 11).Value = 3
$ws.Cells.Item(289, 12).Value = 220
$ws.Cells.Item(289, 13).Value = 660
$ws.Cells.Item(289, 14).Value = "00000000002393"
$ws.Cells.Item(289, 15).Value = "SPW3"

# Row 290
$ws.Cells.Item(290, 8).Value = "K.08.0451"
$ws.Cells.Item(290, 9).Value = "MODULO DE TELEFONE RJ 11 4 FIOS  BRANCA REF 615010BC PIAL PLUS +"
$ws.Cells.Item(290, 10).Value = "UN"
$ws.Cells.Item(290, 11).Value = 20
$ws.Cells.Item(290, 12).Value = 0
$ws.Cells.Item(290, 13).Value = 0

# Row 291
$ws.Cells.Item(291, 8).Value = "K.08.0452"
$ws.Cells.Item(291, 9).Value = "MODULO HDMI  BRNACA REF 615092BC PIAL PLUS +"
$ws.Cells.Item(291, 10).Value = "UN"
$ws.Cells.Item(291, 11).Value = 20
$ws.Cells.Item(291, 12).Value = 0
$ws.Cells.Item(291, 13).Value = 0

# Row 292
$ws.Cells.Item(292, 8).Value = "K.08.0443"
$ws.Cells.Item(292, 9).Value = "MODULO DE TOMADA RJ 45 CAT 6 A LCS 2  PIAL PLUS +  REF 615044BC"
$ws.Cells.Item(292, 10).Value = "UN"
$ws.Cells.Item(292, 11).Value = 10
$ws.Cells.Item(292, 12).Value = 0
$ws.Cells.Item(292, 13).Value = 0

# Row 293
$ws.Cells.Item(293, 8).Value = "K.08.0453"
$ws.Cells.Item(293, 9).Value = "MODULO DE ANTENA DE TV  COAXIAL BRANCA  REF 615030BC PIAL PLUS +"
$ws.Cells.Item(293, 10).Value = "UN"
$ws.Cells.Item(293, 11).Value = 20
$ws.Cells.Item(293, 12).Value = 0
$ws.Cells.Item(293, 13).Value = 0

# Row 294
$ws.Cells.Item(294, 8).Value = "K.08.0826"
$ws.Cells.Item(294, 9).Value = "MÓDULO DE  INTERRUPTOR   ITERMEDIÁRIO   PIAL +  REF 612007 BC"
$ws.Cells.Item(294, 10).Value = "UN"
$ws.Cells.Item(294, 11).Value = 1
$ws.Cells.Item(294, 12).Value = 0
$ws.Cells.Item(294, 13).Value = 0

# Row 295
$ws.Cells.Item(295, 8).Value = "K.08.0907"
$ws.Cells.Item(295, 9).Value = "TOMADA PADRÃO BRASILEIRO 2P+T 20 A - REF.: 0642 19 PIAL AQUATIC"
$ws.Cells.Item(295, 10).Value = "UN"
$ws.Cells.Item(295, 11).Value = 6
$ws.Cells.Item(295, 12).Value = 0
$ws.Cells.Item(295, 13).Value = 0

# Row 296
$ws.Cells.Item(296, 8).Value = "K.08.0908"
$ws.Cells.Item(296, 9).Value = "INTERRUPTOR SIMPLES  20 A - REF.: 642 01 PIAL AQUATIC"
$ws.Cells.Item(296, 10).Value = "UN"
$ws.Cells.Item(296, 11).Value = 1
$ws.Cells.Item(296, 12).Value = 0
$ws.Cells.Item(296, 13).Value = 0

# Row 297
$ws.Cells.Item(297, 8).Value = "K.08.0909"
$ws.Cells.Item(297, 9).Value = "INTERRUPTOR PARALELO  REF.: 642 02  PIAL AQUATIC"
$ws.Cells.Item(297, 10).Value = "UN"
$ws.Cells.Item(297, 11).Value = 6
$ws.Cells.Item(297, 12).Value = 0
$ws.Cells.Item(297, 13).Value = 0

# Row 298
$ws.Cells.Item(298, 8).Value = "O.01.0142"
$ws.Cells.Item(298, 9).Value = "TABUA DE PINUS  1`" X 12`""
$ws.Cells.Item(298, 10).Value = "M"
$ws.Cells.Item(298, 11).Value = 18
$ws.Cells.Item(298, 12).Value = 0
$ws.Cells.Item(298, 13).Value = 0
